# Actualizacion desde MV -datos- : agrega/actualiza filas de indices bursatiles
# (25-10-2021 revisado, y nuevas filas 26-10-2021 .. 02-11-2021)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value (e.g. a "DD-MM-YYYY" string) into a cell
# without Excels automatic text-to-date recognition kicking in. We build the
# text via a formula (so it is unambiguously a string), then copy/paste-values
# it into the destination cell, and finally clear the scratch cell.
function Set-TextValue($cellRef, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = "=""" + $text + """"
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

# --- Update existing row 212 (25-10-2021): revise a few values + append new column data ---
$ws.Range("B212").Value = 35741.2
$ws.Range("C212").Value = 15226.7
$ws.Range("D212").Value = 7222.8
$ws.Range("F212").Value = 6712.9
$ws.Range("G212").Value = 15599.2
$ws.Range("K212").Value = 73876.39999999999
$ws.Range("L212").Value = 4255
$ws.Range("O212").Value = 1492.9
$ws.Range("P212").Value = 89391.8
$ws.Range("Q212").Value = 108714.6
$ws.Range("R212").Value = 1402
$ws.Range("S212").Value = 51833.8
$ws.Range("T212").Value = 21243.1

# --- Append new rows 213-218 (updated data through 02-11-2021) ---
# Row 213: 26-10-2021
Set-TextValue "A213" "26-10-2021"
$ws.Range("B213").Value = 35756.9
$ws.Range("C213").Value = 15235.7
$ws.Range("D213").Value = 7277.6
$ws.Range("E213").Value = 29106
$ws.Range("F213").Value = 6766.5
$ws.Range("G213").Value = 15757.1
$ws.Range("H213").Value = 3049.1
$ws.Range("I213").Value = 4963.1
$ws.Range("J213").Value = 1584.2
$ws.Range("K213").Value = 73909.39999999999
$ws.Range("L213").Value = 4235.9
$ws.Range("M213").Value = 1636
$ws.Range("N213").Value = 17034.3
$ws.Range("O213").Value = 1509.2
$ws.Range("P213").Value = 88907.8
$ws.Range("Q213").Value = 106419.5
$ws.Range("R213").Value = 1413.3
$ws.Range("S213").Value = 52206.6
$ws.Range("T213").Value = 21005.8

# Row 214: 27-10-2021
Set-TextValue "A214" "27-10-2021"
$ws.Range("B214").Value = 35490.7
$ws.Range("C214").Value = 15235.8
$ws.Range("D214").Value = 7253.3
$ws.Range("E214").Value = 29098.2
$ws.Range("F214").Value = 6753.5
$ws.Range("G214").Value = 15705.8
$ws.Range("H214").Value = 3025.5
$ws.Range("I214").Value = 4898.2
$ws.Range("J214").Value = 1583.1
$ws.Range("K214").Value = 73377.89999999999
$ws.Range("L214").Value = 4229.5
$ws.Range("M214").Value = 1627.6
$ws.Range("N214").Value = 17074.6
$ws.Range("O214").Value = 1519.3
$ws.Range("P214").Value = 85353.10000000001
$ws.Range("Q214").Value = 106363.1
$ws.Range("R214").Value = 1404.1
$ws.Range("S214").Value = 51714.6
$ws.Range("T214").Value = 20885.9

# Row 215: 28-10-2021
Set-TextValue "A215" "28-10-2021"
$ws.Range("B215").Value = 35730.5
$ws.Range("C215").Value = 15448.1
$ws.Range("D215").Value = 7249.5
$ws.Range("E215").Value = 28820.1
$ws.Range("F215").Value = 6804.2
$ws.Range("G215").Value = 15696.3
$ws.Range("H215").Value = 3009.6
$ws.Range("I215").Value = 4864.1
$ws.Range("J215").Value = 1566.9
$ws.Range("K215").Value = 73217.39999999999
$ws.Range("L215").Value = 4189.7
$ws.Range("M215").Value = 1624.3
$ws.Range("N215").Value = 17041.6
$ws.Range("O215").Value = 1522
$ws.Range("P215").Value = 86034.2
$ws.Range("Q215").Value = 105705
$ws.Range("R215").Value = 1407.9
$ws.Range("S215").Value = 51248.8
$ws.Range("T215").Value = 20959.5

# Row 216: 29-10-2021
Set-TextValue "A216" "29-10-2021"
$ws.Range("B216").Value = 35819.6
$ws.Range("C216").Value = 15498.4
$ws.Range("D216").Value = 7237.6
$ws.Range("E216").Value = 28892.7
$ws.Range("F216").Value = 6830.3
$ws.Range("G216").Value = 15688.8
$ws.Range("H216").Value = 2970.7
$ws.Range("I216").Value = 4908.8
$ws.Range("J216").Value = 1562.3
$ws.Range("K216").Value = 73586.3
$ws.Range("L216").Value = 4150
$ws.Range("M216").Value = 1623.4
$ws.Range("N216").Value = 16987.4
$ws.Range("P216").Value = 83561
$ws.Range("Q216").Value = 103500.7
$ws.Range("R216").Value = 1394
$ws.Range("S216").Value = 51309.8
$ws.Range("T216").Value = 20737.8

# Row 217: 01-11-2021
Set-TextValue "A217" "01-11-2021"
$ws.Range("B217").Value = 35913.8
$ws.Range("C217").Value = 15595.9
$ws.Range("D217").Value = 7288.6
$ws.Range("E217").Value = 29647.1
$ws.Range("F217").Value = 6893.3
$ws.Range("G217").Value = 15806.3
$ws.Range("H217").Value = 2978.9
$ws.Range("I217").Value = 4890.7
$ws.Range("J217").Value = 1530.9
$ws.Range("L217").Value = 4221.5
$ws.Range("M217").Value = 1613.8
$ws.Range("N217").Value = 17068.2
$ws.Range("O217").Value = 1536.3
$ws.Range("P217").Value = 87460.39999999999
$ws.Range("Q217").Value = 105550.9
$ws.Range("S217").Value = 51653.3

# Row 218: 02-11-2021
Set-TextValue "A218" "02-11-2021"
$ws.Range("D218").Value = 7248.8
$ws.Range("E218").Value = 29520.9
$ws.Range("F218").Value = 6910.8
$ws.Range("G218").Value = 15915.6
$ws.Range("H218").Value = 3013.5
$ws.Range("I218").Value = 4839.9
$ws.Range("J218").Value = 1537.6
$ws.Range("K218").Value = 73753.10000000001
$ws.Range("L218").Value = 4200.6
$ws.Range("M218").Value = 1617.9
$ws.Range("N218").Value = 17066
$ws.Range("O218").Value = 1536.4
